# Update Work Week and Social Spending
# Refresh GDP per Capita data for Zimbabwe (country code 716): revise 1950-2010 figures
# and extend the series through 2016.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing GDP per Capita values for years 1950-2010 (rows 2-62)
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "1117"
$cell.ClearFormats()
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "1151"
$cell.ClearFormats()
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "1154"
$cell.ClearFormats()
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "1211"
$cell.ClearFormats()
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "1231"
$cell.ClearFormats()
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "1288"
$cell.ClearFormats()
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "1422"
$cell.ClearFormats()
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "1473"
$cell.ClearFormats()
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "1444"
$cell.ClearFormats()
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "1474"
$cell.ClearFormats()
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "1495"
$cell.ClearFormats()
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "1524"
$cell.ClearFormats()
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "1497"
$cell.ClearFormats()
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "1436"
$cell.ClearFormats()
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "1519"
$cell.ClearFormats()
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "1568"
$cell.ClearFormats()
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "1541"
$cell.ClearFormats()
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "1618"
$cell.ClearFormats()
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "1592"
$cell.ClearFormats()
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "1731"
$cell.ClearFormats()
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "2043"
$cell.ClearFormats()
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "2157"
$cell.ClearFormats()
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "2268"
$cell.ClearFormats()
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "2283"
$cell.ClearFormats()
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "2275"
$cell.ClearFormats()
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "2235"
$cell.ClearFormats()
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "2163"
$cell.ClearFormats()
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "1946"
$cell.ClearFormats()
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "1964"
$cell.ClearFormats()
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "1930"
$cell.ClearFormats()
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "2064"
$cell.ClearFormats()
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "2243"
$cell.ClearFormats()
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "2240"
$cell.ClearFormats()
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "2190"
$cell.ClearFormats()
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "2067"
$cell.ClearFormats()
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "2128"
$cell.ClearFormats()
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "2107"
$cell.ClearFormats()
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "2004"
$cell.ClearFormats()
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "2114"
$cell.ClearFormats()
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "2181"
$cell.ClearFormats()
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "2160"
$cell.ClearFormats()
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "2222.25533870045"
$cell.ClearFormats()
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "1971.44557700488"
$cell.ClearFormats()
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "1948.91053825787"
$cell.ClearFormats()
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "2110.71198430126"
$cell.ClearFormats()
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "2102.85990985999"
$cell.ClearFormats()
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "2292.4874652967"
$cell.ClearFormats()
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "2326.53888823625"
$cell.ClearFormats()
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "2369.01912391686"
$cell.ClearFormats()
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "2323.00927537187"
$cell.ClearFormats()
$cell = $ws.Range("E52")
$cell.NumberFormat = "@"
$cell.Value = "2211.19619342415"
$cell.ClearFormats()
$cell = $ws.Range("E53")
$cell.NumberFormat = "@"
$cell.Value = "2193.73848050911"
$cell.ClearFormats()
$cell = $ws.Range("E54")
$cell.NumberFormat = "@"
$cell.Value = "2025.31775378879"
$cell.ClearFormats()
$cell = $ws.Range("E55")
$cell.NumberFormat = "@"
$cell.Value = "1700.95585834871"
$cell.ClearFormats()
$cell = $ws.Range("E56")
$cell.NumberFormat = "@"
$cell.Value = "1604.50302473472"
$cell.ClearFormats()
$cell = $ws.Range("E57")
$cell.NumberFormat = "@"
$cell.Value = "1496.03426670055"
$cell.ClearFormats()
$cell = $ws.Range("E58")
$cell.NumberFormat = "@"
$cell.Value = "1455.7286141876"
$cell.ClearFormats()
$cell = $ws.Range("E59")
$cell.NumberFormat = "@"
$cell.Value = "1422.1553208298"
$cell.ClearFormats()
$cell = $ws.Range("E60")
$cell.NumberFormat = "@"
$cell.Value = "1197.52606766824"
$cell.ClearFormats()
$cell = $ws.Range("E61")
$cell.NumberFormat = "@"
$cell.Value = "1285.04658864784"
$cell.ClearFormats()
$cell = $ws.Range("E62")
$cell.NumberFormat = "@"
$cell.Value = "1401.85651907259"
$cell.ClearFormats()

# Add new rows for years 2011-2016 (rows 63-68)
$ws.Range("A63").Value = 716
$ws.Range("B63").Value = "Zimbabwe"
$ws.Range("C63").Value = "GDP per Capita"
$ws.Range("D63").Value = 2011
$cell = $ws.Range("E63")
$cell.NumberFormat = "@"
$cell.Value = "1515"
$cell.ClearFormats()
$ws.Range("A64").Value = 716
$ws.Range("B64").Value = "Zimbabwe"
$ws.Range("C64").Value = "GDP per Capita"
$ws.Range("D64").Value = 2012
$cell = $ws.Range("E64")
$cell.NumberFormat = "@"
$cell.Value = "1604"
$cell.ClearFormats()
$ws.Range("A65").Value = 716
$ws.Range("B65").Value = "Zimbabwe"
$ws.Range("C65").Value = "GDP per Capita"
$ws.Range("D65").Value = 2013
$cell = $ws.Range("E65")
$cell.NumberFormat = "@"
$cell.Value = "1604"
$cell.ClearFormats()
$ws.Range("A66").Value = 716
$ws.Range("B66").Value = "Zimbabwe"
$ws.Range("C66").Value = "GDP per Capita"
$ws.Range("D66").Value = 2014
$cell = $ws.Range("E66")
$cell.NumberFormat = "@"
$cell.Value = "1594"
$cell.ClearFormats()
$ws.Range("A67").Value = 716
$ws.Range("B67").Value = "Zimbabwe"
$ws.Range("C67").Value = "GDP per Capita"
$ws.Range("D67").Value = 2015
$cell = $ws.Range("E67")
$cell.NumberFormat = "@"
$cell.Value = "1560"
$cell.ClearFormats()
$ws.Range("A68").Value = 716
$ws.Range("B68").Value = "Zimbabwe"
$ws.Range("C68").Value = "GDP per Capita"
$ws.Range("D68").Value = 2016
$cell = $ws.Range("E68")
$cell.NumberFormat = "@"
$cell.Value = "1534"
$cell.ClearFormats()
